$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the account summary figures ---
# VALOR MORA total
$ws.Range("E11").Value = 217603
# Cant. Trabajadores (count of workers listed)
$ws.Range("C13").Value = 2

# --- Clear the existing worker rows so stale values/strings are dropped ---
$ws.Range("B16:J22").ClearContents()

# --- Remove the now-blank spacer rows so the table keeps 5 data rows total ---
# (this also shifts the footer rows 27/28 up to 25/26)
$ws.Rows("19:20").Delete()

# --- Re-populate the worker rows with the refreshed database values ---
# Row 16: MARCO ANTONIO LEMUS ARRIETA - periodo 1905
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73187266"
$ws.Range("D16").Value = "MARCO ANTONIO LEMUS ARRIETA"
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 33125
$ws.Range("G16").Value = 781242

# Row 17: MARCO ANTONIO LEMUS ARRIETA - periodo 1906
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73187266"
$ws.Range("D17").Value = "MARCO ANTONIO LEMUS ARRIETA"
$ws.Range("E17").Value = "1906"
$ws.Range("F17").Value = 33125
$ws.Range("G17").Value = 781242

# Row 18: EFRAIN CABALLERO ABELLO - periodo 1906
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "73581077"
$ws.Range("D18").Value = "EFRAIN CABALLERO ABELLO"
$ws.Range("E18").Value = "1906"
$ws.Range("F18").Value = 63200
$ws.Range("G18").Value = 1580000

# Row 19: MARCO ANTONIO LEMUS ARRIETA - periodo 1907
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "73187266"
$ws.Range("D19").Value = "MARCO ANTONIO LEMUS ARRIETA"
$ws.Range("E19").Value = "1907"
$ws.Range("F19").Value = 29166
$ws.Range("G19").Value = 781242

# Row 20: EFRAIN CABALLERO ABELLO - periodo 1907
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "73581077"
$ws.Range("D20").Value = "EFRAIN CABALLERO ABELLO"
$ws.Range("E20").Value = "1907"
$ws.Range("F20").Value = 58987
$ws.Range("G20").Value = 1580000
